$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3, 13, 2, 7),
    @(3, 8, 4, 12),
    @(7, 13, 6, 7),
    @(5, 2, 4, 18),
    @(4, 7, 3, 13),
    @(6, 8, 4, 12),
    @(5, 16, 3, 4),
    @(4, 13, 7, 7),
    @(2, 13, 4, 7),
    @(6, 5, 4, 15),
    @(4, 13, 5, 7),
    @(7, 16, 2, 4),
    @(4, 4, 5, 16),
    @(3, 7, 4, 13),
    @(7, 12, 6, 8),
    @(4, 6, 3, 14),
    @(5, 19, 3, 1)
)

$startRow = 1338
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $rowValues[0]
    $ws.Cells.Item($r, 2).Value2 = $rowValues[1]
    $ws.Cells.Item($r, 3).Value2 = $rowValues[2]
    $ws.Cells.Item($r, 4).Value2 = $rowValues[3]
}

$ws.Application.ActiveWindow.ScrollRow = 1337
$ws.Range("A1355").Select()
